$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Remove now-stale rows 21-23 (sheet shrinks from 23 to 20 rows) ---
$ws.Range("A21:H23").Delete()

# --- Remove all existing hyperlinks up front; we rebuild them below in clean row order ---
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-11-12 06:28:50"
$ws.Range("B2").Value = "専門データ分析:AIコスト最適化設計と厳格な機密保持を必須とするWebシステム開発(段階的継続発注)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5431917"
$ws.Range("G2").Value = 403
$ws.Range("H2").Value = "🔥AI,Ai ◆開発,システム開発"

# Row 3
$ws.Range("A3").Value = "2025-11-12 06:28:50"
$ws.Range("B3").Value = "詳細設計及び、Next.js,node.jsによるWEBアプリケーション開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5427010"
$ws.Range("G3").Value = 245
$ws.Range("H3").Value = "🔥Next.js ◆開発,Node.js ◇アプリ"

# Row 4
$ws.Range("A4").Value = "2025-11-12 06:28:50"
$ws.Range("B4").Value = "<Next.js、バックエンド開発> ガントチャートアプリの改修製造"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5427011"
$ws.Range("G4").Value = 225
$ws.Range("H4").Value = "🔥Next.js ◆開発 ◇アプリ"

# Row 5
$ws.Range("A5").Value = "2025-11-12 06:28:50"
$ws.Range("B5").Value = "【急募】大手保険会社向けスマホアプリ設計書作成依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5431609"
$ws.Range("G5").Value = 98
$ws.Range("H5").Value = "★スマホアプリ ◇アプリ"

# Row 6
$ws.Range("A6").Value = "2025-11-12 06:28:50"
$ws.Range("B6").Value = "ヤフオクの指定出品者の出品物を一括してウォッチリストに登録するツール"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "~ 5,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5431786"
$ws.Range("G6").Value = 65
$ws.Range("H6").Value = "◆ツール"

# Row 7
$ws.Range("A7").Value = "2025-11-12 06:28:50"
$ws.Range("B7").Value = "PHP業務アプリケーションの改修対応"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5426598"
$ws.Range("G7").Value = 58
$ws.Range("H7").Value = "◇アプリ ○PHP"

# Row 8
$ws.Range("A8").Value = "2025-11-12 06:28:50"
$ws.Range("B8").Value = "【案件】既存WordPressサイトの読み込み速度改善"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5432161"
$ws.Range("G8").Value = 50
$ws.Range("H8").Value = "◇サイト ○WordPress"

# Row 9
$ws.Range("A9").Value = "2025-11-12 06:28:50"
$ws.Range("B9").Value = "【急募】ショッピファイでジャケット仕様確定システム構築"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5432465"
$ws.Range("G9").Value = 40
$ws.Range("H9").ClearContents()

# Row 10
$ws.Range("A10").Value = "2025-11-12 06:28:50"
$ws.Range("B10").Value = "小売店向けシステム性能試験"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5430176"
$ws.Range("G10").Value = 40
$ws.Range("H10").ClearContents()

# Row 11
$ws.Range("A11").Value = "2025-11-12 06:28:50"
$ws.Range("B11").Value = "【電気錠制御】オフィスセキュリティシステム刷新の協力者募集"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5431852"
$ws.Range("G11").Value = 33
$ws.Range("H11").ClearContents()

# Row 12
$ws.Range("A12").Value = "2025-11-12 06:28:50"
$ws.Range("B12").Value = "OR(operations research)にて最適化の仕組みの構築 (リモート)"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5427007"
$ws.Range("G12").Value = 25
$ws.Range("H12").ClearContents()

# Row 13
$ws.Range("A13").Value = "2025-11-12 06:28:50"
$ws.Range("B13").Value = "OR(operations research)にて最適化の仕組みの構築(社内常駐)"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5427009"
$ws.Range("G13").Value = 25
$ws.Range("H13").ClearContents()

# Row 14
$ws.Range("A14").Value = "2025-11-12 06:28:50"
$ws.Range("B14").Value = "【急募】楽天市場在庫連動システム(同一店舗内)のエラー修正依頼"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5432212"
$ws.Range("G14").Value = 25
$ws.Range("H14").ClearContents()

# Row 15
$ws.Range("A15").Value = "2025-11-12 06:28:50"
$ws.Range("B15").Value = "【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5431911"
$ws.Range("G15").Value = 25
$ws.Range("H15").ClearContents()

# Row 16
$ws.Range("A16").Value = "2025-11-12 06:28:50"
$ws.Range("B16").Value = "【音楽制作】サイケデリックトランスのトラックを作成してくれる方募集"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5432042"
$ws.Range("G16").Value = 18
$ws.Range("H16").ClearContents()

# Row 17
$ws.Range("A17").Value = "2025-11-12 06:28:50"
$ws.Range("B17").Value = "【スポット案件】HTML途切れ・白画面・Segmentation fault調査対応"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5432323"
$ws.Range("G17").Value = 13
$ws.Range("H17").ClearContents()

# Row 18
$ws.Range("A18").Value = "2025-11-12 06:28:50"
$ws.Range("B18").Value = "初回 MT4用インジケータの修正カスタマイズ(.mq4)"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5432362"
$ws.Range("G18").Value = 10
$ws.Range("H18").ClearContents()

# Row 19
$ws.Range("A19").Value = "2025-11-12 06:28:50"
$ws.Range("B19").Value = "MT4用インジケータの修正カスタマイズ(.mq4)"
$ws.Range("C19").Value = "システム開発"
$ws.Range("D19").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E19").Value = "期限情報なし"
$ws.Range("F19").Value = "https://www.lancers.jp/work/detail/5432305"
$ws.Range("G19").Value = 10
$ws.Range("H19").ClearContents()

# Row 20
$ws.Range("A20").Value = "2025-11-12 06:28:50"
$ws.Range("B20").Value = "【Stable Diffusion】参考動画に沿って約100プロンプト構築"
$ws.Range("C20").Value = "システム開発"
$ws.Range("D20").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E20").Value = "期限情報なし"
$ws.Range("F20").Value = "https://www.lancers.jp/work/detail/5432055"
$ws.Range("G20").Value = 10
$ws.Range("H20").ClearContents()

# --- Rebuild hyperlinks on F2:F20 in row order so relationship IDs stay sequential ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5431917")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5427010")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5427011")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5431609")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5431786")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5426598")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5432161")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5432465")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5430176")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5431852")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5427007")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5427009")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5432212")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5431911")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5432042")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5432323")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5432362")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5432305")
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5432055")

# --- Column H width change (22 -> 27 raw OOXML units); COM ColumnWidth adds ~0.8333 padding ---
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668

$ws.Range("A1").Select()